$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row 1 (B1:C1 -> B1:M1), reusing the existing header style ---
$ws.Range("B1:C1").Copy()
$ws.Range("D1:M1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10
$ws.Range("L1").Value = 11
$ws.Range("M1").Value = 12

# --- Data rows 2-15, columns B..M ---
$data = @(
    @(2, @(87,69,45,30,18,9,58,34,16,1,1,1)),
    @(3, @(57,37,25,19,15,1,1,1,1,11,1,1)),
    @(4, @(115,91,58,37,1,151,118,85,73,37,1,1)),
    @(5, @(82,73,43,54,37,34,1,67,34,1,1,1)),
    @(6, @(25,21,19,11,12,1,15,13,1,16,10,1)),
    @(7, @(12,7,1,27,15,9,1,24,20,9,1,1)),
    @(8, @(56,47,35,25,18,12,5,18,16,7,1,1)),
    @(9, @(17,5,1,1,6,5,1,16,11,1,12,1)),
    @(10, @(56,1,1,1,4,1,11,6,1,36,1,1)),
    @(11, @(33,3,3,303,3,3,3,243,93,3,213,3)),
    @(12, @(79,55,1,1,103,61,43,1,1,61,1,1)),
    @(13, @(55,40,25,56,21,1,1,41,1,1,1,1)),
    @(14, @(1662,342,12,3,3,3,443,3,3,3,3,3)),
    @(15, @(21,1,1,21,1,1,56,1,1,1,16,1))
)

foreach ($entry in $data) {
    $r = $entry[0]
    $vals = $entry[1]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        # column 2 = B
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}

Write-Host "edit applied"
